$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Time/LastLB now measured for just the first layer, Layers 2 -> 1
$ws.Range("L7").Value = 0.52
$ws.Range("M7").Value = 0.52
$ws.Range("O7").Value = 1

# Row 8
$ws.Range("L8").Value = 0.68
$ws.Range("M8").Value = 0.68
$ws.Range("O8").Value = 1

# Row 27: LastLB becomes a formula mirroring Time (L27)
$ws.Range("L27").Value = 0.25
$ws.Range("M27").Formula = "=L27"
$ws.Range("O27").Value = 1

# Row 28
$ws.Range("L28").Value = 0.25
$ws.Range("M28").Formula = "=L28"
$ws.Range("O28").Value = 1

# Row 31
$ws.Range("L31").Value = 0.67
$ws.Range("M31").Formula = "=L31"
$ws.Range("O31").Value = 1

# Row 32
$ws.Range("L32").Value = 1.42
$ws.Range("M32").Formula = "=L32"
$ws.Range("O32").Value = 1

# Row 33
$ws.Range("L33").Value = 3.59
$ws.Range("M33").Formula = "=L33"
$ws.Range("O33").Value = 1

# Update the active selection on the sheet (was E21, now P1)
$ws.Range("P1").Select()
